$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.478.49"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.825.80"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5181"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3864"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08283"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.08%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.381"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.493"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001122"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06633"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.059"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "28.520.42"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.251"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").Value = "2.037.98"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.414"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1096"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.733"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.683"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2231"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02367"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.264"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.766"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6411"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6202"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.798"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.203"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06976"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
